$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark Price (D) column cells as Text first so Excel does not coerce
# numeric-looking strings (e.g. "0.660", "247.69") into floating point numbers,
# which would lose trailing zeros / exact formatting.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = '36.993.43'
$ws.Range("E2").Value = '  -0.24%  '

# Row 3
$ws.Range("D3").Value = '2.047.38'
$ws.Range("E3").Value = '  +0.03%  '

# Row 4
$ws.Range("E4").Value = '  -0.06%  '

# Row 5
$ws.Range("D5").Value = '247.69'
$ws.Range("E5").Value = '  +0.04%  '

# Row 6
$ws.Range("D6").Value = '0.660'
$ws.Range("E6").Value = '  +1.31%  '

# Row 7
$ws.Range("B7").Value = 'Solana'
$ws.Range("C7").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D7").Value = '57.74'
$ws.Range("E7").Value = '  +6.01%  '

# Row 8
$ws.Range("B8").Value = 'USDC'
$ws.Range("C8").Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range("D8").Value = '1.00'
$ws.Range("E8").Value = '  +0.01%  '

# Row 9
$ws.Range("D9").Value = '0.378'
$ws.Range("E9").Value = '  +1.10%  '

# Row 10
$ws.Range("E10").Value = '  +1.33%  '

# Row 11
$ws.Range("E11").Value = '  +1.81%  '

# Row 12
$ws.Range("D12").Value = '15.75'
$ws.Range("E12").Value = '  +5.18%  '

# Row 13
$ws.Range("D13").Value = '2.344.50'
$ws.Range("E13").Value = '  +0.01%  '

# Row 14
$ws.Range("D14").Value = '0.805'
$ws.Range("E14").Value = '  -0.72%  '

# Row 15
$ws.Range("D15").Value = '5.52'
$ws.Range("E15").Value = '  +6.96%  '

# Row 16
$ws.Range("D16").Value = '2.045.78'
$ws.Range("E16").Value = '  -0.04%  '

# Row 17
$ws.Range("D17").Value = '37.053.02'
$ws.Range("E17").Value = '  +0.19%  '

# Row 18
$ws.Range("D18").Value = '16.53'

# Row 19
$ws.Range("D19").Value = '74.21'
$ws.Range("E19").Value = '  +3.56%  '

# Row 20
$ws.Range("D20").Value = '0.0₃0896'
$ws.Range("E20").Value = '  -0.37%  '

# Row 21
$ws.Range("E21").Value = '  +1.55%  '

# Row 22
$ws.Range("D22").Value = '235.18'
$ws.Range("E22").Value = '  -0.26%  '

# Row 23
$ws.Range("E23").Value = '  +0.08%  '

# Row 24
$ws.Range("E24").Value = '  -1.65%  '

# Row 25
$ws.Range("E25").Value = '  +10.38%  '

# Row 26
$ws.Range("D26").Value = '167.15'
$ws.Range("E26").Value = '  -1.28%  '

# Row 27
$ws.Range("D27").Value = '9.09'
$ws.Range("E27").Value = '  +1.42%  '

# Row 28
$ws.Range("E28").Value = '  -1.62%  '

# Row 29
$ws.Range("E29").Value = '  +1.05%  '

# Row 30
$ws.Range("D30").Value = '1.13'
$ws.Range("E30").Value = '  +6.15%  '

# Row 31
$ws.Range("E31").Value = '  +3.16%  '

# Row 32
$ws.Range("E32").Value = '  -1.35%  '

# Row 33
$ws.Range("D33").Value = '4.44'
$ws.Range("E33").Value = '  +2.89%  '

# Row 34
$ws.Range("D34").Value = '0.0883'
$ws.Range("E34").Value = '  +1.44%  '

# Row 35
$ws.Range("E35").Value = '  -0.06%  '

# Row 36
$ws.Range("E36").Value = '  -1.93%  '

# Row 37
$ws.Range("E37").Value = '  -1.42%  '

# Row 38
$ws.Range("B38").Value = 'TrustWalletToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D38").Value = '1.34'
$ws.Range("E38").Value = '  +0.14%  '

# Row 39
$ws.Range("B39").Value = 'Cronos'
$ws.Range("C39").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D39").Value = '0.106'
$ws.Range("E39").Value = '  +2.62%  '

# Row 40
$ws.Range("D40").Value = '3.17'
$ws.Range("E40").Value = '  +14.15%  '

# Row 41
$ws.Range("E41").Value = '  +24.43%  '

# Row 42
$ws.Range("E42").Value = '  -1.37%  '

# Row 43
$ws.Range("D43").Value = '17.23'
$ws.Range("E43").Value = '  -5.54%  '

# Row 44
$ws.Range("E44").Value = '  -0.76%  '

# Row 45
$ws.Range("D45").Value = '95.09'
$ws.Range("E45").Value = '  -0.09%  '

# Row 46
$ws.Range("E46").Value = '  +3.43%  '

# Row 47
$ws.Range("D47").Value = '1.273.92'
$ws.Range("E47").Value = '  -1.13%  '

# Row 48
$ws.Range("E48").Value = '  -1.80%  '

# Row 49
$ws.Range("D49").Value = '2.229.97'
$ws.Range("E49").Value = '  -0.05%  '

# Row 50
$ws.Range("E50").Value = '  -1.34%  '

# Row 51
$ws.Range("D51").Value = '3.45'
$ws.Range("E51").Value = '  -13.95%  '

# Restore original (unstyled/default) formatting on the Price column cells we touched,
# now that the text values are safely stored, so the cell style matches the original workbook.
$ws.Range("D2").ClearFormats()
$ws.Range("D3").ClearFormats()
$ws.Range("D5").ClearFormats()
$ws.Range("D6").ClearFormats()
$ws.Range("D7").ClearFormats()
$ws.Range("D8").ClearFormats()
$ws.Range("D9").ClearFormats()
$ws.Range("D12").ClearFormats()
$ws.Range("D13").ClearFormats()
$ws.Range("D14").ClearFormats()
$ws.Range("D15").ClearFormats()
$ws.Range("D16").ClearFormats()
$ws.Range("D17").ClearFormats()
$ws.Range("D18").ClearFormats()
$ws.Range("D19").ClearFormats()
$ws.Range("D20").ClearFormats()
$ws.Range("D22").ClearFormats()
$ws.Range("D26").ClearFormats()
$ws.Range("D27").ClearFormats()
$ws.Range("D30").ClearFormats()
$ws.Range("D33").ClearFormats()
$ws.Range("D34").ClearFormats()
$ws.Range("D38").ClearFormats()
$ws.Range("D39").ClearFormats()
$ws.Range("D40").ClearFormats()
$ws.Range("D43").ClearFormats()
$ws.Range("D45").ClearFormats()
$ws.Range("D47").ClearFormats()
$ws.Range("D49").ClearFormats()
$ws.Range("D51").ClearFormats()
